$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.838.15"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "2.524.30"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.19"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.81"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.580"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.544"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.02"
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.72"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "2.912.46"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "2.528.27"
$ws.Range("E15").Value = "  -4.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.15"
$ws.Range("E16").Value = "  +6.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.864"
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").Value = "42.832.40"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.98"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.47"
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.40"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.27"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("E25").Value = "  -4.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.86"
$ws.Range("E26").Value = "  -6.52%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.43"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.32"
$ws.Range("E29").Value = "  +10.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.41"
$ws.Range("E30").Value = "  +3.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.01"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.85"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("E33").Value = "  -3.44%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0793"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.09"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("E36").Value = "  -4.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.116"
$ws.Range("E37").Value = "  +2.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.34"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.14"
$ws.Range("E40").Value = "  +3.63%  "
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0302"
$ws.Range("E45").Value = "  -3.06%  "
$ws.Range("D46").Value = "2.037.09"
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.12"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.00"
$ws.Range("E48").Value = "  -2.86%  "
$ws.Range("D49").Value = "2.775.35"
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.190"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.67"
$ws.Range("E51").Value = "  -0.92%  "
